$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2:E25").Value = "fullRNASEQ"

$ws.Range("E24:E25").Select()
$excel.ActiveWindow.ScrollRow = 16
